$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.134.54"
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = "'3.420.56"
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'571.51"
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").Value = "'161.14"
$ws.Range("E6").Value = '  +2.29%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = "'3.421.94"
$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = '  -6.01%  '

$ws.Range("D10").Value = "'7.28"
$ws.Range("E10").Value = '  +1.61%  '

$ws.Range("E11").Value = '  -1.23%  '

$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = '  -2.75%  '

$ws.Range("D13").Value = "'4.013.96"
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").Value = "'27.01"
$ws.Range("E15").Value = '  -1.91%  '

$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = '  -6.37%  '

$ws.Range("D17").Value = "'64.176.14"
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").Value = "'3.429.38"
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").Value = "'6.07"
$ws.Range("E19").Value = '  -3.81%  '

$ws.Range("D20").Value = "'13.55"
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("D21").Value = "'376.35"
$ws.Range("E21").Value = '  -0.53%  '

$ws.Range("D22").Value = "'7.81"
$ws.Range("E22").Value = '  -1.61%  '

$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").Value = "'71.37"
$ws.Range("E24").Value = '  -0.42%  '

$ws.Range("D25").Value = "'0.517"
$ws.Range("E25").Value = '  -4.91%  '

$ws.Range("E26").Value = '  -1.53%  '

$ws.Range("D27").Value = "'9.49"
$ws.Range("E27").Value = '  -3.87%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("D30").Value = "'6.01"
$ws.Range("E30").Value = '  -1.97%  '

$ws.Range("E31").Value = '  -3.57%  '

$ws.Range("E32").Value = '  +0.84%  '

$ws.Range("D33").Value = "'22.83"
$ws.Range("E33").Value = '  -1.28%  '

$ws.Range("D34").Value = "'7.05"
$ws.Range("E34").Value = '  +1.15%  '

$ws.Range("E35").Value = '  -3.51%  '

$ws.Range("E36").Value = '  -0.69%  '

$ws.Range("D37").Value = "'0.854"
$ws.Range("E37").Value = '  +11.03%  '

$ws.Range("E38").Value = '  -3.93%  '

$ws.Range("D39").Value = "'2.810.43"
$ws.Range("E39").Value = '  -2.31%  '

$ws.Range("D40").Value = "'0.0725"
$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = "'42.89"
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = "'25.72"
$ws.Range("E42").Value = '  -1.54%  '

$ws.Range("D43").Value = "'6.46"
$ws.Range("E43").Value = '  -3.04%  '

$ws.Range("D44").Value = "'4.42"
$ws.Range("E44").Value = '  -2.06%  '

$ws.Range("D45").Value = "'25.87"
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").Value = "'0.0304"
$ws.Range("E46").Value = '  -2.96%  '

$ws.Range("E47").Value = '  +9.30%  '

$ws.Range("D48").Value = "'335.03"

$ws.Range("E49").Value = '  -0.93%  '

$ws.Range("D50").Value = "'6.30"
$ws.Range("E50").Value = '  -2.78%  '

$ws.Range("E51").Value = '  -2.68%  '

